$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" (E) and "Valor Mora" (F) values between rows 17 and 18
$e17 = $ws.Range("E17").Value()
$e18 = $ws.Range("E18").Value()
$f17 = $ws.Range("F17").Value()
$f18 = $ws.Range("F18").Value()

$ws.Range("E17").Value = $e18
$ws.Range("E18").Value = $e17
$ws.Range("F17").Value = $f18
$ws.Range("F18").Value = $f17
